$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 20 (shifts old rows 20-40 down to 22-42)
$ws.Rows("20:21").Insert()

# Insert one more new row before (the now-shifted) row 33 (shifts rows 33-42 down to 34-43)
$ws.Rows("33:33").Insert()

# --- Fill in new row 20 ---
$ws.Range("A20").Value = 6
$ws.Range("B20").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44435
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 100114007
$ws.Range("G20").Value = "Jengibre"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 580
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 13000
$ws.Range("M20").Value = 12500
$ws.Range("N20").Value = "$/caja 13 kilos"
$ws.Range("O20").Value = "Perú"
$ws.Range("P20").Value = 962
$ws.Range("Q20").Value = 13
$ws.Range("R20").Value = "Hortaliza"

# --- Fill in new row 21 ---
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44431
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100114007
$ws.Range("G21").Value = "Jengibre"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 260
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 13000
$ws.Range("M21").Value = 12462
$ws.Range("N21").Value = "$/caja 13 kilos"
$ws.Range("O21").Value = "Perú"
$ws.Range("P21").Value = 959
$ws.Range("Q21").Value = 13
$ws.Range("R21").Value = "Hortaliza"

# --- Fill in new row 33 ---
$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44433
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = 100114007
$ws.Range("G33").Value = "Jengibre"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 320
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 13000
$ws.Range("M33").Value = 12531
$ws.Range("N33").Value = "$/caja 13 kilos"
$ws.Range("O33").Value = "Perú"
$ws.Range("P33").Value = 964
$ws.Range("Q33").Value = 13
$ws.Range("R33").Value = "Hortaliza"
